$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: "PASSWORD" header, same format as A1 ("USERNAME")
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "PASSWORD"

# B2: empty value entered with a leading apostrophe (text quote prefix)
$ws.Range("B2").Value = "'"

# Row 1 no longer carries an explicit custom height - autofit restores the default
$ws.Rows(1).AutoFit()

# Update selection to B3 as in the edited workbook
$ws.Range("B3").Select() | Out-Null
